$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook tracks weekly Chirimoya price quotes for "Comercializadora del
# Agro de Limarí" (Coquimbo). Each week contributes a small block of rows
# (Especial/Primera/Segunda/etc. quality grades). This edit represents a new
# weekly refresh: every existing block's data shifts to the next reporting
# week, and a brand-new block (mirroring the oldest week's figures) is
# appended as rows 21-23, extending the sheet from A1:T20 to A1:T23.

$ws.Cells.Item(2,4).Value = 44413
$ws.Cells.Item(2,12).Value = 'Primera'
$ws.Cells.Item(2,13).Value = 200
$ws.Cells.Item(2,14).Value = 2600
$ws.Cells.Item(2,15).Value = 2700
$ws.Cells.Item(2,16).Value = 2650
$ws.Cells.Item(2,17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(2,18).Value = 'Provincia del Elquí'
$ws.Cells.Item(2,19).Value = 2650
$ws.Cells.Item(2,20).Value = 1
$ws.Cells.Item(3,4).Value = 44413
$ws.Cells.Item(3,12).Value = 'Segunda'
$ws.Cells.Item(3,13).Value = 200
$ws.Cells.Item(3,14).Value = 2200
$ws.Cells.Item(3,15).Value = 2300
$ws.Cells.Item(3,16).Value = 2250
$ws.Cells.Item(3,17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(3,18).Value = 'Provincia del Elquí'
$ws.Cells.Item(3,19).Value = 2250
$ws.Cells.Item(3,20).Value = 1
$ws.Cells.Item(4,4).Value = 44454
$ws.Cells.Item(4,12).Value = 'Especial'
$ws.Cells.Item(4,13).Value = 400
$ws.Cells.Item(4,14).Value = 2300
$ws.Cells.Item(4,15).Value = 2400
$ws.Cells.Item(4,16).Value = 2350
$ws.Cells.Item(4,17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(4,19).Value = 2350
$ws.Cells.Item(4,20).Value = 1
$ws.Cells.Item(5,4).Value = 44454
$ws.Cells.Item(5,12).Value = 'Extra (doble especial)'
$ws.Cells.Item(5,14).Value = 2700
$ws.Cells.Item(5,15).Value = 2800
$ws.Cells.Item(5,16).Value = 2750
$ws.Cells.Item(5,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(5,19).Value = 2750
$ws.Cells.Item(6,4).Value = 44454
$ws.Cells.Item(6,12).Value = 'Primera'
$ws.Cells.Item(6,13).Value = 500
$ws.Cells.Item(6,14).Value = 2000
$ws.Cells.Item(6,15).Value = 2100
$ws.Cells.Item(6,16).Value = 2050
$ws.Cells.Item(6,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(6,19).Value = 2050
$ws.Cells.Item(7,4).Value = 44412
$ws.Cells.Item(7,12).Value = 'Primera'
$ws.Cells.Item(7,13).Value = 200
$ws.Cells.Item(7,14).Value = 2600
$ws.Cells.Item(7,15).Value = 2700
$ws.Cells.Item(7,16).Value = 2650
$ws.Cells.Item(7,17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(7,18).Value = 'Provincia del Elquí'
$ws.Cells.Item(7,19).Value = 2650
$ws.Cells.Item(7,20).Value = 1
$ws.Cells.Item(8,4).Value = 44412
$ws.Cells.Item(8,12).Value = 'Segunda'
$ws.Cells.Item(8,13).Value = 240
$ws.Cells.Item(8,14).Value = 2200
$ws.Cells.Item(8,15).Value = 2300
$ws.Cells.Item(8,16).Value = 2250
$ws.Cells.Item(8,17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(8,18).Value = 'Provincia del Elquí'
$ws.Cells.Item(8,19).Value = 2250
$ws.Cells.Item(8,20).Value = 1
$ws.Cells.Item(9,4).Value = 44161
$ws.Cells.Item(9,12).Value = 'Especial'
$ws.Cells.Item(9,13).Value = 240
$ws.Cells.Item(9,14).Value = 13000
$ws.Cells.Item(9,15).Value = 13500
$ws.Cells.Item(9,16).Value = 13250
$ws.Cells.Item(9,19).Value = 1656
$ws.Cells.Item(10,12).Value = 'Primera'
$ws.Cells.Item(10,14).Value = 11000
$ws.Cells.Item(10,15).Value = 11500
$ws.Cells.Item(10,16).Value = 11250
$ws.Cells.Item(10,19).Value = 1406
$ws.Cells.Item(11,12).Value = 'Segunda'
$ws.Cells.Item(11,13).Value = 200
$ws.Cells.Item(11,14).Value = 9000
$ws.Cells.Item(11,15).Value = 9500
$ws.Cells.Item(11,16).Value = 9250
$ws.Cells.Item(11,19).Value = 1156
$ws.Cells.Item(12,4).Value = 44448
$ws.Cells.Item(12,12).Value = 'Especial'
$ws.Cells.Item(12,13).Value = 240
$ws.Cells.Item(12,14).Value = 2400
$ws.Cells.Item(12,15).Value = 2500
$ws.Cells.Item(12,16).Value = 2450
$ws.Cells.Item(12,17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(12,19).Value = 2450
$ws.Cells.Item(12,20).Value = 1
$ws.Cells.Item(13,4).Value = 44448
$ws.Cells.Item(13,12).Value = 'Extra (doble especial)'
$ws.Cells.Item(13,13).Value = 240
$ws.Cells.Item(13,14).Value = 2700
$ws.Cells.Item(13,15).Value = 2800
$ws.Cells.Item(13,16).Value = 2750
$ws.Cells.Item(13,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(13,19).Value = 2750
$ws.Cells.Item(14,4).Value = 44448
$ws.Cells.Item(14,12).Value = 'Primera'
$ws.Cells.Item(14,13).Value = 600
$ws.Cells.Item(14,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(15,4).Value = 44167
$ws.Cells.Item(15,13).Value = 360
$ws.Cells.Item(15,14).Value = 12500
$ws.Cells.Item(15,15).Value = 13000
$ws.Cells.Item(15,16).Value = 12750
$ws.Cells.Item(15,17).Value = '$/bandeja 8 kilos'
$ws.Cells.Item(15,19).Value = 1594
$ws.Cells.Item(15,20).Value = 8
$ws.Cells.Item(16,4).Value = 44167
$ws.Cells.Item(16,12).Value = 'Primera'
$ws.Cells.Item(16,13).Value = 300
$ws.Cells.Item(16,14).Value = 10500
$ws.Cells.Item(16,15).Value = 11000
$ws.Cells.Item(16,16).Value = 10750
$ws.Cells.Item(16,17).Value = '$/bandeja 8 kilos'
$ws.Cells.Item(16,19).Value = 1344
$ws.Cells.Item(16,20).Value = 8
$ws.Cells.Item(17,4).Value = 44167
$ws.Cells.Item(17,12).Value = 'Segunda'
$ws.Cells.Item(17,13).Value = 240
$ws.Cells.Item(17,14).Value = 8500
$ws.Cells.Item(17,15).Value = 9000
$ws.Cells.Item(17,16).Value = 8750
$ws.Cells.Item(17,17).Value = '$/bandeja 8 kilos'
$ws.Cells.Item(17,19).Value = 1094
$ws.Cells.Item(17,20).Value = 8
$ws.Cells.Item(18,4).Value = 44168
$ws.Cells.Item(18,13).Value = 300
$ws.Cells.Item(18,14).Value = 13000
$ws.Cells.Item(18,15).Value = 13500
$ws.Cells.Item(18,16).Value = 13250
$ws.Cells.Item(18,19).Value = 1656
$ws.Cells.Item(19,4).Value = 44168
$ws.Cells.Item(19,14).Value = 11000
$ws.Cells.Item(19,15).Value = 11500
$ws.Cells.Item(19,16).Value = 11250
$ws.Cells.Item(19,19).Value = 1406
$ws.Cells.Item(20,4).Value = 44168
$ws.Cells.Item(20,13).Value = 200

# --- Append new rows 21-23 (shifted-in data) ---
# Row 21
$ws.Cells.Item(21,1).Value = 2
$ws.Cells.Item(21,2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(21,3).Value = 'Coquimbo'
$ws.Cells.Item(21,4).NumberFormat = $ws.Cells.Item(2,4).NumberFormat()
$ws.Cells.Item(21,4).Value = 44160
$ws.Cells.Item(21,5).Value = 4
$ws.Cells.Item(21,6).Value = 'Fruta'
$ws.Cells.Item(21,7).Value = 100107
$ws.Cells.Item(21,8).Value = 'Otros'
$ws.Cells.Item(21,9).Value = 100107002
$ws.Cells.Item(21,10).Value = 'Chirimoya'
$ws.Cells.Item(21,11).Value = 'Cultivar IV Región'
$ws.Cells.Item(21,12).Value = 'Especial'
$ws.Cells.Item(21,13).Value = 240
$ws.Cells.Item(21,14).Value = 12500
$ws.Cells.Item(21,15).Value = 13000
$ws.Cells.Item(21,16).Value = 12750
$ws.Cells.Item(21,17).Value = '$/bandeja 8 kilos'
$ws.Cells.Item(21,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(21,19).Value = 1594
$ws.Cells.Item(21,20).Value = 8
# Row 22
$ws.Cells.Item(22,1).Value = 2
$ws.Cells.Item(22,2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(22,3).Value = 'Coquimbo'
$ws.Cells.Item(22,4).NumberFormat = $ws.Cells.Item(2,4).NumberFormat()
$ws.Cells.Item(22,4).Value = 44160
$ws.Cells.Item(22,5).Value = 4
$ws.Cells.Item(22,6).Value = 'Fruta'
$ws.Cells.Item(22,7).Value = 100107
$ws.Cells.Item(22,8).Value = 'Otros'
$ws.Cells.Item(22,9).Value = 100107002
$ws.Cells.Item(22,10).Value = 'Chirimoya'
$ws.Cells.Item(22,11).Value = 'Cultivar IV Región'
$ws.Cells.Item(22,12).Value = 'Primera'
$ws.Cells.Item(22,13).Value = 300
$ws.Cells.Item(22,14).Value = 10500
$ws.Cells.Item(22,15).Value = 11000
$ws.Cells.Item(22,16).Value = 10750
$ws.Cells.Item(22,17).Value = '$/bandeja 8 kilos'
$ws.Cells.Item(22,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(22,19).Value = 1344
$ws.Cells.Item(22,20).Value = 8
# Row 23
$ws.Cells.Item(23,1).Value = 2
$ws.Cells.Item(23,2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(23,3).Value = 'Coquimbo'
$ws.Cells.Item(23,4).NumberFormat = $ws.Cells.Item(2,4).NumberFormat()
$ws.Cells.Item(23,4).Value = 44160
$ws.Cells.Item(23,5).Value = 4
$ws.Cells.Item(23,6).Value = 'Fruta'
$ws.Cells.Item(23,7).Value = 100107
$ws.Cells.Item(23,8).Value = 'Otros'
$ws.Cells.Item(23,9).Value = 100107002
$ws.Cells.Item(23,10).Value = 'Chirimoya'
$ws.Cells.Item(23,11).Value = 'Cultivar IV Región'
$ws.Cells.Item(23,12).Value = 'Segunda'
$ws.Cells.Item(23,13).Value = 240
$ws.Cells.Item(23,14).Value = 8500
$ws.Cells.Item(23,15).Value = 9000
$ws.Cells.Item(23,16).Value = 8750
$ws.Cells.Item(23,17).Value = '$/bandeja 8 kilos'
$ws.Cells.Item(23,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(23,19).Value = 1094
$ws.Cells.Item(23,20).Value = 8

# Dimension / used range will be recalculated automatically by Excel based on
# the populated cells (A1:T23).
